$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.856696666666667
$ws.Range("H2").Value = 8.57009
$ws.Range("I2").Value = 0.05747862151401942
$ws.Range("J2").Value = 0.05747862151401942
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 455.5963751580089
$ws.Range("R2").Value = 4100.36737642208
$ws.Range("S2").Value = 0.01714751142635162
$ws.Range("T2").Value = 0.01714751142635162

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.856696666666667
$ws.Range("H3").Value = 8.57009
$ws.Range("I3").Value = 0.05747862151401942
$ws.Range("J3").Value = 0.05747862151401942
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 492.94728318491
$ws.Range("R3").Value = 4436.52554866419
$ws.Range("S3").Value = 0.01855330646138403
$ws.Range("T3").Value = 0.01855330646138403

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.856696666666667
$ws.Range("H4").Value = 8.57009
$ws.Range("I4").Value = 0.05747862151401942
$ws.Range("J4").Value = 0.05747862151401942
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 212.5031108189478
$ws.Range("R4").Value = 1912.52799737053
$ws.Range("S4").Value = 0.007998087165727343
$ws.Range("T4").Value = 0.007998087165727344

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.856696666666667
$ws.Range("H5").Value = 8.57009
$ws.Range("I5").Value = 0.05747862151401942
$ws.Range("J5").Value = 0.05747862151401942
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 166.8728340506522
$ws.Range("R5").Value = 1501.85550645587
$ws.Range("S5").Value = 0.006280677337783543
$ws.Range("T5").Value = 0.006280677337783544

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.856696666666667
$ws.Range("H6").Value = 8.57009
$ws.Range("I6").Value = 0.05747862151401942
$ws.Range("J6").Value = 0.05747862151401942
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 199.2437827598133
$ws.Range("R6").Value = 1793.19404483832
$ws.Range("S6").Value = 0.007499039122772879
$ws.Range("T6").Value = 0.00749903912277288

# Row 7
$ws.Range("I7").Value = 0.7708435061432634
$ws.Range("J7").Value = 0.7708435061432632
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 6109.984859802226
$ws.Range("R7").Value = 54989.86373822004
$ws.Range("S7").Value = 0.229964593467096
$ws.Range("T7").Value = 0.2299645934670959

# Row 8
$ws.Range("I8").Value = 0.7708435061432634
$ws.Range("J8").Value = 0.7708435061432632
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.2488176547476082
$ws.Range("T8").Value = 0.2488176547476082

# Row 9
$ws.Range("I9").Value = 0.7708435061432634
$ws.Range("J9").Value = 0.7708435061432632
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 2849.870764038323
$ws.Range("R9").Value = 25648.83687634491
$ws.Range("S9").Value = 0.1072620287486357
$ws.Range("T9").Value = 0.1072620287486357

# Row 10
$ws.Range("I10").Value = 0.7708435061432634
$ws.Range("J10").Value = 0.7708435061432632
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 2237.924937853516
$ws.Range("R10").Value = 20141.32444068165
$ws.Range("S10").Value = 0.08422991387903674
$ws.Range("T10").Value = 0.08422991387903674

# Row 11
$ws.Range("I11").Value = 0.7708435061432634
$ws.Range("J11").Value = 0.7708435061432632
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 2672.050443004459
$ws.Range("R11").Value = 24048.45398704013
$ws.Range("S11").Value = 0.1005693153008867
$ws.Range("T11").Value = 0.1005693153008867

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 2.081608
$ws.Range("H12").Value = 6.244823999999999
$ws.Range("I12").Value = 0.04188332620983732
$ws.Range("J12").Value = 0.04188332620983732
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 331.9824153421653
$ws.Range("R12").Value = 2987.841738079488
$ws.Range("S12").Value = 0.01249499023879035
$ws.Range("T12").Value = 0.01249499023879035

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 2.081608
$ws.Range("H13").Value = 6.244823999999999
$ws.Range("I13").Value = 0.04188332620983732
$ws.Range("J13").Value = 0.04188332620983732
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 359.1991478231759
$ws.Range("R13").Value = 3232.792330408583
$ws.Range("S13").Value = 0.01351936017817853
$ws.Range("T13").Value = 0.01351936017817853

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 2.081608
$ws.Range("H14").Value = 6.244823999999999
$ws.Range("I14").Value = 0.04188332620983732
$ws.Range("J14").Value = 0.04188332620983732
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 154.8460432173786
$ws.Range("R14").Value = 1393.614388956408
$ws.Range("S14").Value = 0.005828018922394756
$ws.Range("T14").Value = 0.005828018922394757

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 2.081608
$ws.Range("H15").Value = 6.244823999999999
$ws.Range("I15").Value = 0.04188332620983732
$ws.Range("J15").Value = 0.04188332620983732
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 121.5963285131813
$ws.Range("R15").Value = 1094.366956618632
$ws.Range("S15").Value = 0.00457658257675786
$ws.Range("T15").Value = 0.004576582576757861

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 2.081608
$ws.Range("H16").Value = 6.244823999999999
$ws.Range("I16").Value = 0.04188332620983732
$ws.Range("J16").Value = 0.04188332620983732
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 145.184281195328
$ws.Range("R16").Value = 1306.658530757952
$ws.Range("S16").Value = 0.005464374293715819
$ws.Range("T16").Value = 0.005464374293715821

# Row 17
$ws.Range("G17").Value = 4.573220666666667
$ws.Range("H17").Value = 13.719662
$ws.Range("I17").Value = 0.09201621679565497
$ws.Range("J17").Value = 0.09201621679565496
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 729.3538662479716
$ws.Range("R17").Value = 6564.184796231744
$ws.Range("S17").Value = 0.02745106071356101
$ws.Range("T17").Value = 0.02745106071356101

# Row 18
$ws.Range("G18").Value = 4.573220666666667
$ws.Range("H18").Value = 13.719662
$ws.Range("I18").Value = 0.09201621679565497
$ws.Range("J18").Value = 0.09201621679565496
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 789.148084689338
$ws.Range("R18").Value = 7102.332762204042
$ws.Range("S18").Value = 0.02970156598502524
$ws.Range("T18").Value = 0.02970156598502524

# Row 19
$ws.Range("G19").Value = 4.573220666666667
$ws.Range("H19").Value = 13.719662
$ws.Range("I19").Value = 0.09201621679565497
$ws.Range("J19").Value = 0.09201621679565496
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 340.1913929007171
$ws.Range("R19").Value = 3061.722536106454
$ws.Range("S19").Value = 0.01280395568311618
$ws.Range("T19").Value = 0.01280395568311618

# Row 20
$ws.Range("G20").Value = 4.573220666666667
$ws.Range("H20").Value = 13.719662
$ws.Range("I20").Value = 0.09201621679565497
$ws.Range("J20").Value = 0.09201621679565496
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 267.1429215045629
$ws.Range("R20").Value = 2404.286293541066
$ws.Range("S20").Value = 0.01005459338296915
$ws.Range("T20").Value = 0.01005459338296915

# Row 21
$ws.Range("G21").Value = 4.573220666666667
$ws.Range("H21").Value = 13.719662
$ws.Range("I21").Value = 0.09201621679565497
$ws.Range("J21").Value = 0.09201621679565496
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 318.9648364329973
$ws.Range("R21").Value = 2870.683527896976
$ws.Range("S21").Value = 0.01200504103098339
$ws.Range("T21").Value = 0.01200504103098339

# Row 22
$ws.Range("G22").Value = 1.877589
$ws.Range("H22").Value = 5.632767
$ws.Range("I22").Value = 0.037778329337225
$ws.Range("J22").Value = 0.037778329337225
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 299.444723137056
$ws.Range("R22").Value = 2695.002508233504
$ws.Range("S22").Value = 0.01127035264442688
$ws.Range("T22").Value = 0.01127035264442687

# Row 23
$ws.Range("G23").Value = 1.877589
$ws.Range("H23").Value = 5.632767
$ws.Range("I23").Value = 0.037778329337225
$ws.Range("J23").Value = 0.037778329337225
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 323.993935823733
$ws.Range("R23").Value = 2915.945422413597
$ws.Range("S23").Value = 0.01219432379083192
$ws.Range("T23").Value = 0.01219432379083192

# Row 24
$ws.Range("G24").Value = 1.877589
$ws.Range("H24").Value = 5.632767
$ws.Range("I24").Value = 0.037778329337225
$ws.Range("J24").Value = 0.037778329337225
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 139.669537894971
$ws.Range("R24").Value = 1257.025841054739
$ws.Range("S24").Value = 0.005256813108174186
$ws.Range("T24").Value = 0.005256813108174186

# Row 25
$ws.Range("G25").Value = 1.877589
$ws.Range("H25").Value = 5.632767
$ws.Range("I25").Value = 0.037778329337225
$ws.Range("J25").Value = 0.037778329337225
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 109.678637311509
$ws.Range("R25").Value = 987.107735803581
$ws.Range("S25").Value = 0.004128030399437463
$ws.Range("T25").Value = 0.004128030399437462

# Row 26
$ws.Range("G26").Value = 1.877589
$ws.Range("H26").Value = 5.632767
$ws.Range("I26").Value = 0.037778329337225
$ws.Range("J26").Value = 0.037778329337225
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 130.954727953224
$ws.Range("R26").Value = 1178.592551579016
$ws.Range("S26").Value = 0.004928809394354554
$ws.Range("T26").Value = 0.004928809394354554

